$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---- Sheet 展览 (sheet1): bump '想去人数' (F column) counts ----
$ws1.Range("F2").Value = 2397
$ws1.Range("F3").Value = 613
$ws1.Range("F4").Value = 219
$ws1.Range("F5").Value = 379
$ws1.Range("F6").Value = 637
$ws1.Range("F8").Value = 843
$ws1.Range("F10").Value = 872
$ws1.Range("F11").Value = 398
$ws1.Range("F13").Value = 416
$ws1.Range("F16").Value = 22486
$ws1.Range("F17").Value = 1121
$ws1.Range("F18").Value = 110
$ws1.Range("F19").Value = 308
$ws1.Range("F22").Value = 199
$ws1.Range("F24").Value = 27
$ws1.Range("F25").Value = 290
$ws1.Range("F26").Value = 21
$ws1.Range("F27").Value = 385
$ws1.Range("F28").Value = 170

# ---- Sheet 演出 (sheet2): bump '想去人数' (F column) counts ----
$ws2.Range("F8").Value = 3492
$ws2.Range("F10").Value = 132
$ws2.Range("F14").Value = 127
$ws2.Range("F16").Value = 4073

# ---- Sheet 本地生活 (sheet3): bump '想去人数' (F column) counts ----
$ws3.Range("F4").Value = 672

# ---- Sheet 全部类型 (sheet4): bump '想去人数' (F column) counts on rows not otherwise rewritten ----
$ws4.Range("F4").Value = 38
$ws4.Range("F5").Value = 2397
$ws4.Range("F6").Value = 672
$ws4.Range("F17").Value = 843
$ws4.Range("F19").Value = 872
$ws4.Range("F20").Value = 398
$ws4.Range("F22").Value = 416
$ws4.Range("F25").Value = 22486
$ws4.Range("F27").Value = 3492
$ws4.Range("F29").Value = 132
$ws4.Range("F31").Value = 1121
$ws4.Range("F32").Value = 110
$ws4.Range("F33").Value = 308
$ws4.Range("F38").Value = 199
$ws4.Range("F40").Value = 27
$ws4.Range("F41").Value = 127
$ws4.Range("F43").Value = 290
$ws4.Range("F44").Value = 21
$ws4.Range("F45").Value = 385
$ws4.Range("F46").Value = 170
$ws4.Range("F47").Value = 4073

# ---- Sheet 全部类型 (sheet4): rows 7-10 content cascades up by one
# (row 7's original '砂糖桔动漫荟STJ02' entry is removed; later rows shift up,
#  and a new '奥斯卡·罗曼耶卓' concert entry lands on row 10) ----
# Row 7
$ws4.Range("C7").Value = '广州·第六届淋唔到动漫嘉年华'
$ws4.Range("D7").Value = '沿江东三路15号 广州1978文化创意园'
$ws4.Range("F7").Value = 219
$ws4.Range("G7").Value = 50
$ws4.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=85554'
$ws4.Range("I7").Value = '//i2.hdslb.com/bfs/openplatform/202405/MtLwpx7j1715570717678.jpeg'

# Row 8
$ws4.Range("C8").Value = '广州·第十一届樱漫动漫嘉年华'
$ws4.Range("D8").Value = '奥体南路12号 优托邦(奥体旗舰店)'
$ws4.Range("E8").Value = '2024.06.22 10:00-06.22 17:00'
$ws4.Range("F8").Value = 379
$ws4.Range("G8").Value = 54
$ws4.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=86075'
$ws4.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202405/bTLocQ1C1716258923709.jpeg'

# Row 9
$ws4.Range("B9").Value = "'2024-06-23"
$ws4.Range("C9").Value = '广州·第五人格ONLY'
$ws4.Range("E9").Value = '2024.06.23 10:00-06.23 17:00'
$ws4.Range("F9").Value = 637
$ws4.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=86276'
$ws4.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202406/sk6wpxN91717486689960.jpeg'

# Row 10
$ws4.Range("B10").Value = "'2024-06-28"
$ws4.Range("C10").Value = '广州·奥斯卡·罗曼耶卓（O叔）钢琴独奏音乐会'
$ws4.Range("D10").Value = '晴波路33号 广州星海音乐厅'
$ws4.Range("E10").Value = '2024.06.28 20:00-06.28 21:30'
$ws4.Range("F10").Value = 185
$ws4.Range("G10").Value = 480
$ws4.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=84545'
$ws4.Range("I10").Value = '//i2.hdslb.com/bfs/openplatform/202404/XK8EYxGv1712890578712.jpeg'

